$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018")
if (-not $ws) { $ws = $wb.ActiveSheet }

$ws.Range("A45").Value = "Added Henry Moore Page"
$ws.Range("B45").Value = 43506
$ws.Range("C45").Value = 0.75

# Match the date number formatting already used by the cell above (B44)
$ws.Range("B44").Copy()
$ws.Range("B45").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$wb.Save()
